$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correct household composition: HH type 2 consists of 2x person type 1
$ws.Range("D5").Value = 2
$ws.Range("D6").Value = 0

# Update the active selection to D10
$ws.Range("D10").Select()
